$wb = $excel.ActiveWorkbook

# Unhide "Population Definitions" sheet
$popDefSheet = $wb.Worksheets.Item("Population Definitions")
$popDefSheet.Visible = $true

# Delete the "Metadata" sheet entirely
$excel.DisplayAlerts = $false
$metaSheet = $wb.Worksheets.Item("Metadata")
$metaSheet.Delete() | Out-Null
$excel.DisplayAlerts = $true

# Update formulas referencing Population Definitions B2 -> A2 on Parameters and State Variables sheets
$paramSheet = $wb.Worksheets.Item("Parameters")
$stateSheet = $wb.Worksheets.Item("State Variables")

foreach ($ws in @($paramSheet, $stateSheet)) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Columns.Item(1).Cells) {
        if ($cell.HasFormula) {
            $f = $cell.Formula
            if ($f -like "*Population Definitions*`$B`$2*") {
                $cell.Formula = $f.Replace("`$B`$2", "`$A`$2")
            }
        }
    }
}

# Update the "State Variables" sheet's zoom level
$stateSheet.Activate()
$excel.ActiveWindow.Zoom = 85

# Re-point the selection/active sheet to match the edited workbook state
$paramSheet.Activate()
$paramSheet.Range("L19").Select() | Out-Null
